# Auto-generated Excel COM-interop script
# Applies the "scheduled runner" data refresh to the Phoenix_Profits workbook:
# updates current market price / profit columns (H..N) across several leve rows
# on multiple crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1819.0588
$ws.Range("J70").Value = 3598
$ws.Range("L70").Value = 10794
$ws.Range("N70").Value = -11334

$ws.Range("H73").Value = 1819.0588
$ws.Range("J73").Value = 3598
$ws.Range("L73").Value = 10794
$ws.Range("N73").Value = -12666

$ws.Range("H100").Value = 2658.6924
$ws.Range("I100").Value = 3063.125
$ws.Range("K100").Value = 3063.125
$ws.Range("M100").Value = -2522.125

$ws.Range("H132").Value = 2212.6099
$ws.Range("I132").Value = 1803.1621
$ws.Range("K132").Value = 5409.4863
$ws.Range("M132").Value = -2879.4863

$ws.Range("H137").Value = 2362.8572
$ws.Range("I137").Value = 1678.8485
$ws.Range("K137").Value = 5036.5455
$ws.Range("M137").Value = -2486.5455

$ws.Range("H139").Value = 78075
$ws.Range("J139").Value = 78075
$ws.Range("L139").Value = 78075
$ws.Range("N139").Value = -88355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 14015.1
$ws.Range("I31").Value = 11270.125
$ws.Range("K31").Value = 11270.125
$ws.Range("M31").Value = -10976.125

$ws.Range("H40").Value = 34939.234
$ws.Range("J40").Value = 34748.25
$ws.Range("L40").Value = 34748.25
$ws.Range("N40").Value = -35100.25

$ws.Range("H42").Value = 22000
$ws.Range("J42").Value = 22000
$ws.Range("L42").Value = 22000
$ws.Range("N42").Value = -22972

$ws.Range("H63").Value = 3245.4
$ws.Range("I63").Value = 1493.8667
$ws.Range("K63").Value = 1493.8667
$ws.Range("M63").Value = -807.8667

$ws.Range("H66").Value = 3245.4
$ws.Range("I66").Value = 1493.8667
$ws.Range("K66").Value = 7469.333500000001
$ws.Range("M66").Value = -4037.333500000001

$ws.Range("H132").Value = 3925.875
$ws.Range("I132").Value = 3925.875
$ws.Range("K132").Value = 11777.625
$ws.Range("M132").Value = -9247.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 57881.5
$ws.Range("I28").Value = 49900
$ws.Range("J28").Value = 60542
$ws.Range("K28").Value = 49900
$ws.Range("L28").Value = 60542
$ws.Range("M28").Value = -49606
$ws.Range("N28").Value = -61130

$ws.Range("H86").Value = 35498.54
$ws.Range("J86").Value = 69991.62
$ws.Range("L86").Value = 69991.62
$ws.Range("N86").Value = -72237.62

$ws.Range("H89").Value = 35498.54
$ws.Range("J89").Value = 69991.62
$ws.Range("L89").Value = 349958.1
$ws.Range("N89").Value = -361190.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 48444
$ws.Range("J28").Value = 48444
$ws.Range("L28").Value = 48444
$ws.Range("N28").Value = -48934

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1706.8667
$ws.Range("J12").Value = 1846.3077
$ws.Range("L12").Value = 5538.9231
$ws.Range("N12").Value = -5884.9231

$ws.Range("H17").Value = 133.33333
$ws.Range("I17").Value = 133.33333
$ws.Range("K17").Value = 399.99999
$ws.Range("M17").Value = -230.99999

$ws.Range("H68").Value = 5638.8
$ws.Range("I68").Value = 486.125
$ws.Range("J68").Value = 26249.5
$ws.Range("K68").Value = 1458.375
$ws.Range("L68").Value = 78748.5
$ws.Range("M68").Value = -647.375
$ws.Range("N68").Value = -80370.5

$ws.Range("H69").Value = 4948.6665
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14189

$ws.Range("H71").Value = 5638.8
$ws.Range("I71").Value = 486.125
$ws.Range("J71").Value = 26249.5
$ws.Range("K71").Value = 4375.125
$ws.Range("L71").Value = 236245.5
$ws.Range("M71").Value = -319.125
$ws.Range("N71").Value = -244357.5

$ws.Range("H72").Value = 4948.6665
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40944

$ws.Range("H75").Value = 1824.875
$ws.Range("I75").Value = 2320
$ws.Range("J75").Value = 999.6667
$ws.Range("K75").Value = 6960
$ws.Range("L75").Value = 2999.0001
$ws.Range("M75").Value = -5962
$ws.Range("N75").Value = -4995.0001

$ws.Range("H78").Value = 1824.875
$ws.Range("I78").Value = 2320
$ws.Range("J78").Value = 999.6667
$ws.Range("K78").Value = 20880
$ws.Range("L78").Value = 8997.0003
$ws.Range("M78").Value = -15888
$ws.Range("N78").Value = -18981.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5978.8
$ws.Range("I3").Value = 949.5
$ws.Range("J3").Value = 9331.666999999999
$ws.Range("K3").Value = 949.5
$ws.Range("L3").Value = 9331.666999999999
$ws.Range("M3").Value = -833.5
$ws.Range("N3").Value = -9563.666999999999

$ws.Range("H10").Value = 50015000
$ws.Range("I10").Value = 100000000
$ws.Range("J10").Value = 29999
$ws.Range("K10").Value = 100000000
$ws.Range("L10").Value = 29999
$ws.Range("M10").Value = -99999831
$ws.Range("N10").Value = -30337

$ws.Range("H70").Value = 4693.778
$ws.Range("I70").Value = 4755
$ws.Range("J70").Value = 4632.5557
$ws.Range("K70").Value = 4755
$ws.Range("L70").Value = 4632.5557
$ws.Range("M70").Value = -4485
$ws.Range("N70").Value = -5172.5557

$ws.Range("H73").Value = 4693.778
$ws.Range("I73").Value = 4755
$ws.Range("J73").Value = 4632.5557
$ws.Range("K73").Value = 4755
$ws.Range("L73").Value = 4632.5557
$ws.Range("M73").Value = -3819
$ws.Range("N73").Value = -6504.5557

$ws.Range("H80").Value = 12359.6
$ws.Range("I80").Value = 4200
$ws.Range("J80").Value = 14399.5
$ws.Range("K80").Value = 4200
$ws.Range("L80").Value = 14399.5
$ws.Range("M80").Value = -3202
$ws.Range("N80").Value = -16395.5

$ws.Range("H83").Value = 12359.6
$ws.Range("I83").Value = 4200
$ws.Range("J83").Value = 14399.5
$ws.Range("K83").Value = 21000
$ws.Range("L83").Value = 71997.5
$ws.Range("M83").Value = -16008
$ws.Range("N83").Value = -81981.5

$ws.Range("H109").Value = 25075.691
$ws.Range("J109").Value = 25075.691
$ws.Range("L109").Value = 25075.691
$ws.Range("N109").Value = -27155.691

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3182.617
$ws.Range("I132").Value = 2833.4146
$ws.Range("K132").Value = 8500.2438
$ws.Range("M132").Value = -5970.2438

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 404
$ws.Range("I7").Value = 404
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 404
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -291
$ws.Range("N7").ClearContents()

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H96").Value = 2170.7334
$ws.Range("I96").Value = 2488.3
$ws.Range("J96").Value = 1535.6
$ws.Range("K96").Value = 2488.3
$ws.Range("L96").Value = 1535.6
$ws.Range("M96").Value = -1115.3
$ws.Range("N96").Value = -4281.6

$ws.Range("H133").Value = 79999.5
$ws.Range("J133").Value = 79999.5
$ws.Range("L133").Value = 79999.5
$ws.Range("N133").Value = -90119.5
